# Order of entry matches the original authoring order so that the
# generated shared-strings table lines up with the source workbook.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New "Link" column values (column D), entered in the same order the
# author originally typed them (D5, D7, D3, D6, D2). D4 already had a
# link (Seeed XIAO BLE) and is left unchanged.
$ws.Range("D5").Value = "https://es.aliexpress.com/item/33029465106.html"
$ws.Range("D7").Value = "https://es.aliexpress.com/item/4000917776872.html"
$ws.Range("D3").Value = "https://es.aliexpress.com/item/1005001308084552.html"
$ws.Range("D6").Value = "https://es.aliexpress.com/item/4000685483225.html"
$ws.Range("D2").Value = "https://es.aliexpress.com/item/32959996455.html"

# New column E notes
$ws.Range("E2").Value = "1 pack = 10u"
$ws.Range("E3").Value = "do not solder"
$ws.Range("E4").Value = "do not solder"
$ws.Range("E5").Value = "do not solder"
$ws.Range("E6").Value = "do not solder"
$ws.Range("E7").Value = "do not solder"

$ws.Range("F7").Select()

$wb.Save()
